$d = $word.ActiveDocument

$replacements = @(
    @('39×54=2106', '46×60=2760'),
    @('43×15=645', '70×93=6510'),
    @('36×85=3060', '43×80=3440'),
    @('68×53=3604', '88×90=7920'),
    @('42×89=3738', '89×66=5874'),
    @('80×61=4880', '94×40=3760'),
    @('96×65=6240', '67×99=6633'),
    @('53×72=3816', '32×21=672'),
    @('95×75=7125', '13×14=182'),
    @('52×52=2704', '99×47=4653'),
    @('91×82=7462', '81×99=8019'),
    @('74×25=1850', '83×11=913'),
    @('21×87=1827', '33×49=1617'),
    @('23×59=1357', '79×19=1501'),
    @('86×71=6106', '57×25=1425'),
    @('95×28=2660', '80×92=7360'),
    @('97×40=3880', '94×67=6298'),
    @('15×58=870', '72×33=2376'),
    @('78×83=6474', '91×98=8918'),
    @('28×12=336', '69×69=4761'),
    @('75×17=1275', '20×61=1220'),
    @('74×70=5180', '26×61=1586'),
    @('87×65=5655', '51×76=3876'),
    @('51×64=3264', '73×26=1898'),
    @('34×21=714', '42×51=2142'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

